$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F37").Value = 40
$ws.Range("G37").Value = 2801.6
$ws.Range("F50").Value = 32
$ws.Range("G50").Value = 2993.28
$ws.Range("F55").Value = 20
$ws.Range("G55").Value = 706.4
$ws.Range("F56").Value = 30
$ws.Range("G56").Value = 296.1
$ws.Range("B61").Value = 25544.86
$ws.Range("F109").Value = 74
$ws.Range("G109").Value = 4693.82
$ws.Range("F111").Value = 243
$ws.Range("G111").Value = 15479.1
$ws.Range("F116").Value = 152
$ws.Range("G116").Value = 2514.08
$ws.Range("B133").Value = 207171.03
$ws.Range("F187").Value = 3
$ws.Range("G187").Value = 6649.5
$ws.Range("B195").Value = 42946.86
$ws.Range("F253").Value = 30
$ws.Range("G253").Value = 5739.9
$ws.Range("B259").Value = 18682.31
$ws.Range("B314").Value = 57077
$ws.Range("D314").Value = 93.08
$ws.Range("E314").Value = 111.2
$ws.Range("F314").Value = 1
$ws.Range("G314").Value = 93.08
$ws.Range("B315").Value = 61610
$ws.Range("D315").Value = 102.71
$ws.Range("E315").Value = 122.71
$ws.Range("F315").Value = 91
$ws.Range("G315").Value = 9346.610000000001
$ws.Range("F321").Value = 18
$ws.Range("G321").Value = 2171.7
$ws.Range("F326").Value = 29
$ws.Range("G326").Value = 5566.84
$ws.Range("F328").Value = 123
$ws.Range("G328").Value = 17746.44
$ws.Range("F332").Value = 206
$ws.Range("G332").Value = 9813.84
$ws.Range("F354").Value = 37
$ws.Range("G354").Value = 3740.7
$ws.Range("B380").Value = 255647.19
$ws.Range("F385").Value = 8
$ws.Range("G385").Value = 1626.48
$ws.Range("B389").Value = 23013.46
$ws.Range("F440").Value = 125
$ws.Range("G440").Value = 3717.5
$ws.Range("B447").Value = 38385.58
$ws.Range("F491").Value = 475
$ws.Range("G491").Value = 6388.75
$ws.Range("F494").Value = 264
$ws.Range("G494").Value = 6943.2
$ws.Range("F496").Value = 304
$ws.Range("G496").Value = 4994.72
$ws.Range("F497").Value = 281
$ws.Range("G497").Value = 3599.61
$ws.Range("F501").Value = 89
$ws.Range("G501").Value = 1731.94
$ws.Range("F506").Value = 324
$ws.Range("G506").Value = 8521.200000000001
$ws.Range("F507").Value = 254
$ws.Range("G507").Value = 4173.22
$ws.Range("F508").Value = 547
$ws.Range("G508").Value = 8057.31
$ws.Range("B509").Value = 94614.07000000001
$ws.Range("F559").Value = 321
$ws.Range("G559").Value = 2150.7
$ws.Range("F560").Value = 296
$ws.Range("G560").Value = 4892.88
$ws.Range("B563").Value = 36784.66
$ws.Range("F572").Value = 66
$ws.Range("G572").Value = 4085.4
$ws.Range("B584").Value = 23248.12
$ws.Range("F622").Value = 23
$ws.Range("G622").Value = 1155.29
$ws.Range("B640").Value = 208500.32
$ws.Range("F642").Value = 105
$ws.Range("G642").Value = 13707.75
$ws.Range("F643").Value = 69
$ws.Range("G643").Value = 12283.38
$ws.Range("F646").Value = 9
$ws.Range("G646").Value = 244.8
$ws.Range("B649").Value = 53481.23
$ws.Range("F669").Value = 114
$ws.Range("G669").Value = 1808.04
$ws.Range("F671").Value = 147
$ws.Range("G671").Value = 6347.46
$ws.Range("F673").Value = 62
$ws.Range("G673").Value = 2677.16
$ws.Range("F674").Value = 16
$ws.Range("G674").Value = 529.76
$ws.Range("F675").Value = 154
$ws.Range("G675").Value = 6649.72
$ws.Range("B677").Value = 20727.16
$ws.Range("F680").Value = 12
$ws.Range("G680").Value = 906.72
$ws.Range("F682").Value = 18
$ws.Range("G682").Value = 1476.72
$ws.Range("F683").Value = 27
$ws.Range("G683").Value = 2448.36
$ws.Range("F685").Value = 36
$ws.Range("G685").Value = 2253.96
$ws.Range("F687").Value = 34
$ws.Range("G687").Value = 3963.72
$ws.Range("F689").Value = 29
$ws.Range("G689").Value = 2755
$ws.Range("F693").Value = 13
$ws.Range("G693").Value = 1291.03
$ws.Range("F694").Value = 16
$ws.Range("G694").Value = 1761.76
$ws.Range("B695").Value = 44241.35
$ws.Range("F710").Value = 3
$ws.Range("G710").Value = 7186.38
$ws.Range("B716").Value = 103725.52
$ws.Range("F722").Value = 7
$ws.Range("G722").Value = 598.5
$ws.Range("B724").Value = 598.5
$ws.Range("F755").Value = 228
$ws.Range("G755").Value = 18595.68
$ws.Range("F756").Value = 125
$ws.Range("G756").Value = 5982.5
$ws.Range("F757").Value = 25
$ws.Range("G757").Value = 2039
$ws.Range("F758").Value = 254
$ws.Range("G758").Value = 33147
$ws.Range("F761").Value = 30
$ws.Range("G761").Value = 3346.2
$ws.Range("F763").Value = 108
$ws.Range("G763").Value = 2345.76
$ws.Range("F764").Value = 366
$ws.Range("G764").Value = 13644.48
$ws.Range("F771").Value = 491
$ws.Range("G771").Value = 66289.91
$ws.Range("F772").Value = 20
$ws.Range("G772").Value = 748.4
$ws.Range("F773").Value = 567
$ws.Range("G773").Value = 68442.57000000001
$ws.Range("B775").Value = 248554.08
$ws.Range("F800").Value = 10
$ws.Range("G800").Value = 374
$ws.Range("B801").Value = 567.86
$ws.Range("F852").Value = 645
$ws.Range("G852").Value = 19498.35
$ws.Range("F853").Value = 3158
$ws.Range("G853").Value = 515101.38
$ws.Range("F857").Value = 170
$ws.Range("G857").Value = 13113.8
$ws.Range("F860").Value = 128
$ws.Range("G860").Value = 16458.24
$ws.Range("B861").Value = 631995.91
$ws.Range("B867").Value = 3507816.98
$ws.Range("B868").Value = 3507816.98
